$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Change the management-related skill for Rafael Ferreira from
# "Ability to Influence People" to "Scheduling"
$ws.Range("F8").Value = "Scheduling"

# Reflect the selection move recorded in the saved file
$ws.Range("F8").Select()
